$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.047.52'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.061.47'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.92'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.670'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.69'
$ws.Range("E7").Value = '  +5.39%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.30'
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.386'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  +7.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.108'
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '16.45'
$ws.Range("E13").Value = '  +9.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.368.97'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.808'
$ws.Range("E15").Value = '  -2.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.61'
$ws.Range("E16").Value = '  +8.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.064.52'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.045.50'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.53'
$ws.Range("E19").Value = '  +14.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.76'
$ws.Range("E20").Value = '  +3.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0921'
$ws.Range("E21").Value = '  +8.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.46'
$ws.Range("E22").Value = '  +4.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.46'
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.29'
$ws.Range("E26").Value = '  +12.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.30'
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.31'
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.28'
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  +3.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.77'
$ws.Range("E32").Value = '  +5.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0620'
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +7.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0898'
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.28'
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.113'
$ws.Range("E39").Value = '  +18.88%  '
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.84'
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.70'
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.69'
$ws.Range("E46").Value = '  +16.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.07'
$ws.Range("E47").Value = '  -22.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.48'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.293.35'
$ws.Range("E49").Value = '  -2.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.92'
$ws.Range("E51").Value = '  -1.15%  '
